$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (prices, percentages) stay text,
# matching the source workbook which stores them as inline strings.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '70.642.90'
$ws.Range('E2').Value = '  +7.22%  '
$ws.Range('D3').Value = '3.626.63'
$ws.Range('E3').Value = '  +7.05%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').Value = '592.57'
$ws.Range('E5').Value = '  +5.02%  '
$ws.Range('D6').Value = '191.03'
$ws.Range('E6').Value = '  +8.22%  '
$ws.Range('D7').Value = '0.642'
$ws.Range('E7').Value = '  +1.67%  '
$ws.Range('D8').Value = '3.610.48'
$ws.Range('E8').Value = '  +6.88%  '
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').Value = '0.179'
$ws.Range('E10').Value = '  +2.85%  '
$ws.Range('D11').Value = '0.661'
$ws.Range('E11').Value = '  +4.22%  '
$ws.Range('D12').Value = '57.67'
$ws.Range('E12').Value = '  +7.11%  '
$ws.Range('D13').Value = '0.0000293'
$ws.Range('E13').Value = '  +5.31%  '
$ws.Range('D14').Value = '9.71'
$ws.Range('E14').Value = '  +5.20%  '
$ws.Range('D15').Value = '4.213.91'
$ws.Range('E15').Value = '  +7.30%  '
$ws.Range('D16').Value = '3.633.83'
$ws.Range('E16').Value = '  +7.31%  '
$ws.Range('D17').Value = '19.34'
$ws.Range('E17').Value = '  +6.14%  '
$ws.Range('D18').Value = '70.630.75'
$ws.Range('E18').Value = '  +7.24%  '
$ws.Range('D19').Value = '12.56'
$ws.Range('E19').Value = '  +5.69%  '
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('E21').Value = '  +5.24%  '
$ws.Range('D22').Value = '491.12'
$ws.Range('E22').Value = '  +5.55%  '
$ws.Range('E23').Value = '  +11.48%  '
$ws.Range('D24').Value = '16.78'
$ws.Range('E24').Value = '  +15.71%  '
$ws.Range('E25').Value = '  +8.81%  '
$ws.Range('D26').Value = '90.51'
$ws.Range('E26').Value = '  +0.91%  '
$ws.Range('D27').Value = '3.09'
$ws.Range('E27').Value = '  +5.48%  '
$ws.Range('D28').Value = '11.21'
$ws.Range('E28').Value = '  +5.47%  '
$ws.Range('D29').Value = '9.35'
$ws.Range('E29').Value = '  +7.31%  '
$ws.Range('D30').Value = '32.21'
$ws.Range('E30').Value = '  +3.38%  '
$ws.Range('D31').Value = '7.70'
$ws.Range('E31').Value = '  +16.68%  '
$ws.Range('D32').Value = '12.22'
$ws.Range('E32').Value = '  +6.56%  '
$ws.Range('D33').Value = '615.71'
$ws.Range('E33').Value = '  +5.96%  '
$ws.Range('E34').Value = '  +7.71%  '
$ws.Range('D35').Value = '65.33'
$ws.Range('E35').Value = '  +4.89%  '
$ws.Range('D36').Value = '0.0₃0827'
$ws.Range('E36').Value = '  +10.92%  '
$ws.Range('D37').Value = '0.148'
$ws.Range('E37').Value = '  +3.70%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').Value = '0.403'
$ws.Range('E39').Value = '  +6.76%  '
$ws.Range('D40').Value = '37.86'
$ws.Range('E40').Value = '  +5.08%  '
$ws.Range('E41').Value = '  +1.37%  '
$ws.Range('D42').Value = '3.353.01'
$ws.Range('E42').Value = '  +7.98%  '
$ws.Range('D43').Value = '3.07'
$ws.Range('E43').Value = '  +8.02%  '
$ws.Range('D44').Value = '0.0444'
$ws.Range('E44').Value = '  +6.49%  '
$ws.Range('E45').Value = '  +8.73%  '
$ws.Range('E46').Value = '  +6.25%  '
$ws.Range('D47').Value = '0.137'
$ws.Range('E47').Value = '  +2.28%  '
$ws.Range('D48').Value = '9.18'
$ws.Range('E48').Value = '  +7.87%  '
$ws.Range('E49').Value = '  +7.06%  '
$ws.Range('D50').Value = '3.33'
$ws.Range('E50').Value = '  +5.74%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').Value = '143.08'
$ws.Range('E51').Value = '  +1.84%  '
